$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, matching the style of the other header cells
# (copy format first, then overwrite the value so the text isn't clobbered)
$ws.Range("G1:G1").Copy($ws.Range("H1:H1"))
$ws.Range("H1").Value = "Label"

# Updated (refit) prediction/error values for rows 2-11 (100 iterations block)
$ws.Range("D2").Value = 0.4505016930177695
$ws.Range("E2").Value = 0.4505016930177695

$ws.Range("D4").Value = 0.5384472988088268
$ws.Range("E4").Value = 0.5384472988088268

$ws.Range("D5").Value = 0.313668279149673
$ws.Range("E5").Value = 0.313668279149673

$ws.Range("D6").Value = 0.5638485017813335
$ws.Range("E6").Value = 0.5638485017813335

$ws.Range("D7").Value = 0.4980513574126038
$ws.Range("E7").Value = 0.5019486425873962

$ws.Range("D9").Value = 0.5652496355381388
$ws.Range("E9").Value = 0.4347503644618612

$ws.Range("D10").Value = 0.5499176644591492
$ws.Range("E10").Value = 0.4500823355408508

$ws.Range("D11").Value = 0.3974988161800964
$ws.Range("E11").Value = 0.6025011838199037
$ws.Range("F11").Value = 0.6460192799568176

# New "Label" column values for each data row (0 = Control, 1 = MDD)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
